$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-15 22:18:52"
$ws.Range("E3").Value = "2026-02-15 22:18:55"
$ws.Range("I3").Value = "2.8 mm"
$ws.Range("O3").Value = "-4.8 °C"
$ws.Range("E4").Value = "2026-02-15 22:18:58"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "73%"
$ws.Range("E5").Value = "2026-02-15 22:19:01"
$ws.Range("I5").Value = "8.6 mm"
$ws.Range("E6").Value = "2026-02-15 22:19:04"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "62%"
$ws.Range("E7").Value = "2026-02-15 22:19:07"
$ws.Range("E8").Value = "2026-02-15 22:19:09"
$ws.Range("E9").Value = "2026-02-15 22:19:12"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "54%"
$ws.Range("O9").Value = "10.7 °C"
$ws.Range("E10").Value = "2026-02-15 22:19:15"
$ws.Range("O10").Value = "7.3 °C"
$ws.Range("E11").Value = "2026-02-15 22:19:18"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "47%"
$ws.Range("O11").Value = "6.9 °C"
$ws.Range("E12").Value = "2026-02-15 22:19:21"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "60%"
$ws.Range("E13").Value = "2026-02-15 22:19:24"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "40%"
$ws.Range("E14").Value = "2026-02-15 22:19:27"
$ws.Range("O14").Value = "10.8 °C"
$ws.Range("E15").Value = "2026-02-15 22:19:30"
$ws.Range("O15").Value = "10.4 °C"
$ws.Range("E16").Value = "2026-02-15 22:19:32"
$ws.Range("E17").Value = "2026-02-15 22:19:35"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "41%"
$ws.Range("O17").Value = "3.1 °C"
$ws.Range("E18").Value = "2026-02-15 22:19:38"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "73%"
$ws.Range("O18").Value = "7.5 °C"
$ws.Range("E19").Value = "2026-02-15 22:19:41"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "74%"
$ws.Range("O19").Value = "3.6 °C"
$ws.Range("E20").Value = "2026-02-15 22:19:44"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "63%"
$ws.Range("O20").Value = "-2.5 °C"
$ws.Range("E21").Value = "2026-02-15 22:19:46"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "41%"
$ws.Range("E22").Value = "2026-02-15 22:19:49"
$ws.Range("E23").Value = "2026-02-15 22:19:52"
$ws.Range("I23").Value = "5.5 mm"
$ws.Range("E24").Value = "2026-02-15 22:19:55"
$ws.Range("E25").Value = "2026-02-15 22:19:57"
$ws.Range("E26").Value = "2026-02-15 22:20:00"
$ws.Range("E27").Value = "2026-02-15 22:20:03"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "51%"
$ws.Range("E28").Value = "2026-02-15 22:20:05"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "60%"
$ws.Range("E29").Value = "2026-02-15 22:20:08"
$ws.Range("E30").Value = "2026-02-15 22:20:10"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "58%"
$ws.Range("E31").Value = "2026-02-15 22:20:13"
$ws.Range("O31").Value = "10.2 °C"
$ws.Range("E32").Value = "2026-02-15 22:20:16"
$ws.Range("O32").Value = "4.0 °C"
$ws.Range("E33").Value = "2026-02-15 22:20:18"
$ws.Range("E34").Value = "2026-02-15 22:20:21"
$ws.Range("O34").Value = "1.4 °C"
$ws.Range("E35").Value = "2026-02-15 22:20:24"
$ws.Range("O35").Value = "4.3 °C"
$ws.Range("E36").Value = "2026-02-15 22:20:27"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "51%"
$ws.Range("E37").Value = "2026-02-15 22:20:30"
$ws.Range("O37").Value = "5.8 °C"
$ws.Range("E38").Value = "2026-02-15 22:20:33"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "69%"
$ws.Range("E39").Value = "2026-02-15 22:20:36"
$ws.Range("E40").Value = "2026-02-15 22:20:39"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "42%"
$ws.Range("J40").Value = "1016.3 hPa"
$ws.Range("O40").Value = "8.5 °C"
$ws.Range("E41").Value = "2026-02-15 22:20:41"
$ws.Range("O41").Value = "12.7 °C"
$ws.Range("E42").Value = "2026-02-15 22:20:44"
$ws.Range("E43").Value = "2026-02-15 22:20:47"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "69%"
$ws.Range("E44").Value = "2026-02-15 22:20:50"
$ws.Range("I44").Value = "5.9 mm"
$ws.Range("O44").Value = "-3.7 °C"
$ws.Range("E45").Value = "2026-02-15 22:20:53"
$ws.Range("I45").Value = "4.4 mm"
$ws.Range("J45").Value = "1023.2 hPa"
$ws.Range("E46").Value = "2026-02-15 22:20:56"
